$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 ("Wait") to make room for the new
# "Idle" animation entry at the top of the table.
$ws.Range("C4:F4").EntireRow.Insert()

# The insert leaves the new row 4 with default/blank formatting; copy the
# formatting (borders/fill/alignment/number format) from the row right below
# (the old row 4, now shifted to row 5, which still carries the correct data
# row style) so the new row matches the rest of the table.
$ws.Range("C5:F5").Copy()
$ws.Range("C4:F4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "Idle" row.
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = "Idle"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 36

# Existing rows, renumbered / re-ranged to follow the new "Idle" entry.
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Wait"
$ws.Range("E5").Value = 37
$ws.Range("F5").Value = 125

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "Greet"
$ws.Range("E6").Value = 126
$ws.Range("F6").Value = 188

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Talk"
$ws.Range("E7").Value = 189
$ws.Range("F7").Value = 317

# Match the saved selection state recorded in the authored workbook.
$ws.Range("F7").Select()
